$p = $ppt.ActivePresentation

# --- Refresh the cached "update automatically" date field text that PowerPoint
# rewrites (on the slide master and every slide layout) whenever the deck is
# opened and saved again on a later day. ---
$m = $p.SlideMaster

for ($mi = 1; $mi -le $m.Shapes.Count; $mi++) {
    $mShape = $m.Shapes.Item($mi)
    if ($mShape.HasTextFrame) {
        $mRange = $mShape.TextFrame.TextRange
        if ($mRange.Text -eq "2024/4/12") {
            $mRange.Text = "2024/5/7"
        }
    }
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $lShape = $layout.Shapes.Item($si)
        if ($lShape.HasTextFrame) {
            $lRange = $lShape.TextFrame.TextRange
            if ($lRange.Text -eq "2024/4/12") {
                $lRange.Text = "2024/5/7"
            }
        }
    }
}

# --- Add the new third slide (blank layout, same as slides 1-2) with the new
# note about designing a generic GC architecture. ---
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 12)

$newBox = $newSlide.Shapes.AddTextbox(1, 25.474645669291338, 43.20314960629921, 589.671811023622, 50.892204724409446)
$newBox.Name = "TextBox 1"
$newBox.TextFrame.WordWrap = -1
$newBox.TextFrame.AutoSize = 1
$newBox.TextFrame.TextRange.Text = "设计通用的GC架构，问题是：inuse变量放在哪里？是插入到用户结构体里还是另开一个结构体？如果另开的话，递归遍历子孙接口如何设计？"
$newBox.Height = 50.892204724409446

Write-Output "done"
